# Updates the "Student loans" column (column 3) of the regression-results
# table with refreshed model estimates, and corrects the "N" figure in the
# same column (per commit: "continuing to work on the r&r").
#
# Each target value lives in its own table cell, so cells are addressed
# directly (Table.Cell(row, col), 1-indexed) and only the visible text of
# the cell is replaced -- this avoids accidentally matching a numeric
# substring that also appears inside a neighboring confidence-interval
# cell elsewhere in the document (e.g. "-0.006" also occurs inside
# "[-0.006, 0.163]").

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $newText) {
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters so only the
    # visible run text is replaced; formatting of the run is preserved.
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $newText
}

# row, col (1-indexed; row 1 is the header row, col 3 is "Student loans")
Set-CellText 2  3 "0.060"
Set-CellText 3  3 "[-0.002, 0.122]"
Set-CellText 4  3 "0.978"
Set-CellText 5  3 "[0.775, 1.182]"
Set-CellText 6  3 "-0.236"
Set-CellText 7  3 "[-0.287, -0.184]"
Set-CellText 8  3 "-0.002"
Set-CellText 9  3 "[-0.061, 0.058]"
Set-CellText 10 3 "-0.070"
Set-CellText 11 3 "[-0.095, -0.044]"
Set-CellText 12 3 "0.166"
Set-CellText 13 3 "[0.006, 0.326]"
Set-CellText 14 3 "-0.749"
Set-CellText 15 3 "[-0.838, -0.661]"
Set-CellText 16 3 "7.068"
Set-CellText 17 3 "[6.631, 7.505]"
Set-CellText 18 3 "1757"
Set-CellText 19 3 "0.46"
Set-CellText 20 3 "1.68"

Write-Host "Done applying edits."
